$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2371795.28
$ws.Range("C7").Value = -46.61819540548615
$ws.Range("D7").Value = 2421
$ws.Range("E7").Value = 2421
$ws.Range("F7").Value = 979.6758694754233
$ws.Range("G7").Value = 4.426363717314175
